$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: extend with P1=14, Q1=15 (style matches existing header cells)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2: update existing values B2:O2 and add P2:Q2
$row2 = @(2,2,2,1,1,1,2,2,2,1,2,2,2,1,2,2)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = 2 + $i  # B = column 2
    $ws.Cells.Item(2, $col).Value = $row2[$i]
}

# Rows 3-25: update I,K,M columns and add P,Q columns
for ($r = 3; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P new = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q new = 2
}
